$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '37.088.71'
$ws.Range("E2").Value = '  +1.25%  '
$ws.Range("D3").Value = '2.058.40'
$ws.Range("E3").Value = '  -2.36%  '
$ws.Range("E4").Value = '  +0.24%  '
$c = $ws.Range("D5")
$c.Value = '''248.61'
$c.ClearFormats()
$ws.Range("E5").Value = '  -1.72%  '
$ws.Range("E6").Value = '  -0.90%  '
$ws.Range("E7").Value = '  -0.03%  '
$c = $ws.Range("D8")
$c.Value = '''55.48'
$c.ClearFormats()
$ws.Range("E8").Value = '  +15.93%  '
$c = $ws.Range("D9")
$c.Value = '''61.38'
$c.ClearFormats()
$ws.Range("E9").Value = '  +3.14%  '
$ws.Range("E10").Value = '  +1.68%  '
$c = $ws.Range("D11")
$c.Value = '''0.0799'
$c.ClearFormats()
$ws.Range("E11").Value = '  +6.91%  '
$ws.Range("E12").Value = '  +5.49%  '
$c = $ws.Range("D13")
$c.Value = '''15.12'
$c.ClearFormats()
$ws.Range("E13").Value = '  +5.63%  '
$ws.Range("D14").Value = '2.358.13'
$ws.Range("E14").Value = '  -2.33%  '
$c = $ws.Range("D15")
$c.Value = '''0.815'
$c.ClearFormats()
$ws.Range("E15").Value = '  -1.74%  '
$ws.Range("E16").Value = '  +2.71%  '
$ws.Range("D17").Value = '2.064.50'
$ws.Range("E17").Value = '  -1.98%  '
$ws.Range("D18").Value = '37.043.15'
$ws.Range("E18").Value = '  +1.07%  '
$ws.Range("D19").Value = '0.0₃0940'
$ws.Range("E19").Value = '  +13.00%  '
$c = $ws.Range("D20")
$c.Value = '''72.33'
$c.ClearFormats()
$ws.Range("E20").Value = '  -1.23%  '
$c = $ws.Range("D21")
$c.Value = '''14.20'
$c.ClearFormats()
$ws.Range("E21").Value = '  +6.64%  '
$ws.Range("E22").Value = '  +4.18%  '
$c = $ws.Range("D23")
$c.Value = '''237.03'
$c.ClearFormats()
$ws.Range("E23").Value = '  -1.56%  '
$ws.Range("E24").Value = '  +0.07%  '
$ws.Range("E25").Value = '  -1.08%  '
$c = $ws.Range("D26")
$c.Value = '''170.56'
$c.ClearFormats()
$ws.Range("E26").Value = '  -0.73%  '
$c = $ws.Range("D27")
$c.Value = '''9.06'
$c.ClearFormats()
$ws.Range("E27").Value = '  -1.76%  '
$c = $ws.Range("D28")
$c.Value = '''20.18'
$c.ClearFormats()
$ws.Range("E28").Value = '  -6.29%  '
$ws.Range("E29").Value = '  -0.81%  '
$ws.Range("E30").Value = '  +0.00%  '
$c = $ws.Range("D31")
$c.Value = '''4.56'
$c.ClearFormats()
$ws.Range("E31").Value = '  +2.08%  '
$c = $ws.Range("D32")
$c.Value = '''1.05'
$c.ClearFormats()
$ws.Range("E32").Value = '  +11.35%  '
$ws.Range("E33").Value = '  +3.19%  '
$c = $ws.Range("D34")
$c.Value = '''4.36'
$c.ClearFormats()
$ws.Range("E34").Value = '  +7.10%  '
$ws.Range("E35").Value = '  +0.12%  '
$c = $ws.Range("D36")
$c.Value = '''0.0856'
$c.ClearFormats()
$ws.Range("E36").Value = '  -4.80%  '
$ws.Range("E37").Value = '  -3.34%  '
$ws.Range("E38").Value = '  -6.58%  '
$ws.Range("E39").Value = '  +1.27%  '
$c = $ws.Range("D40")
$c.Value = '''0.105'
$c.ClearFormats()
$ws.Range("E40").Value = '  +24.42%  '
$c = $ws.Range("D41")
$c.Value = '''17.94'
$c.ClearFormats()
$ws.Range("E41").Value = '  +11.05%  '
$ws.Range("E42").Value = '  -0.10%  '
$ws.Range("E43").Value = '  -3.40%  '
$ws.Range("B44").Value = 'Aave'
$ws.Range("C44").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$c = $ws.Range("D44")
$c.Value = '''96.38'
$c.ClearFormats()
$ws.Range("E44").Value = '  -1.79%  '
$ws.Range("B45").Value = 'FTXToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$c = $ws.Range("D45")
$c.Value = '''4.34'
$c.ClearFormats()
$ws.Range("E45").Value = '  +47.88%  '
$c = $ws.Range("D46")
$c.Value = '''2.77'
$c.ClearFormats()
$ws.Range("E46").Value = '  +0.50%  '
$ws.Range("B47").Value = 'RenderToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$c = $ws.Range("D47")
$c.Value = '''2.42'
$c.ClearFormats()
$ws.Range("E47").Value = '  +6.34%  '
$ws.Range("B48").Value = 'Gas'
$ws.Range("C48").Value = 'https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas'
$c = $ws.Range("D48")
$c.Value = '''13.30'
$c.ClearFormats()
$ws.Range("E48").Value = '  -53.42%  '
$ws.Range("D49").Value = '1.297.32'
$ws.Range("E49").Value = '  -3.22%  '
$ws.Range("E50").Value = '  +2.90%  '
$ws.Range("B51").Value = 'THORChain'
$ws.Range("C51").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$c = $ws.Range("D51")
$c.Value = '''4.04'
$c.ClearFormats()
$ws.Range("E51").Value = '  +4.73%  '
